$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.646.21"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "3.677.64"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "662.66"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.424"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "3.673.66"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.60%  "
$ws.Range("D15").Value = "4.360.16"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("E16").Value = "  +4.65%  "
$ws.Range("D17").Value = "96.409.15"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +11.23%  "
$ws.Range("D19").Value = "3.657.33"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("E22").Value = "  -3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "531.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "3.875.16"
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("E30").Value = "  +13.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("E35").Value = "  +13.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.594"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.86%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "635.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "44.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +34.07%  "
$ws.Range("E42").Value = "  +5.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.968"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.472"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +26.80%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.59%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0458"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.15%  "
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.64%  "
